$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "trainingaudio/25_tapapi1.wav"
$ws.Range("B2").Value = "pngimages/25_apple.png"
$ws.Range("C2").Value = "trainingaudio/13_kopopi1.wav"
$ws.Range("D2").Value = "pngimages/13_toast.png"
$ws.Range("E2").Value = 0.5
$ws.Range("F2").Value = -0.5

$ws.Range("A3").Value = "trainingaudio/11_tokiko1.wav"
$ws.Range("B3").Value = "pngimages/11_compass.png"
$ws.Range("C3").Value = "trainingaudio/24_takopa1.wav"
$ws.Range("D3").Value = "pngimages/24_banana.png"
$ws.Range("E3").Value = 0.5
$ws.Range("F3").Value = -0.5

$ws.Range("A4").Value = "trainingaudio/09_tipata2.wav"
$ws.Range("B4").Value = "pngimages/09_plane.png"
$ws.Range("C4").Value = "trainingaudio/10_tokiti1.wav"
$ws.Range("D4").Value = "pngimages/10_backpack.png"
$ws.Range("E4").Value = -0.5
$ws.Range("F4").Value = 0.5

$ws.Range("A5").Value = "trainingaudio/17_kotako2.wav"
$ws.Range("B5").Value = "pngimages/17_cracker.png"
$ws.Range("C5").Value = "trainingaudio/08_tipako2.wav"
$ws.Range("D5").Value = "pngimages/08_bell.png"
$ws.Range("E5").Value = -0.5
$ws.Range("F5").Value = 0.5

$ws.Range("A6").Value = "trainingaudio/03_kikita3.wav"
$ws.Range("B6").Value = "pngimages/03_box.png"
$ws.Range("C6").Value = "trainingaudio/09_tipata2.wav"
$ws.Range("D6").Value = "pngimages/09_plane.png"
$ws.Range("E6").Value = -0.5
$ws.Range("F6").Value = 0.5

$ws.Range("A7").Value = "trainingaudio/17_kotako2.wav"
$ws.Range("B7").Value = "pngimages/17_cracker.png"
$ws.Range("C7").Value = "trainingaudio/25_tapapi1.wav"
$ws.Range("D7").Value = "pngimages/25_apple.png"
$ws.Range("E7").Value = -0.5
$ws.Range("F7").Value = 0.5

$ws.Range("A8").Value = "trainingaudio/03_kikita3.wav"
$ws.Range("B8").Value = "pngimages/03_box.png"
$ws.Range("C8").Value = "trainingaudio/24_takopa1.wav"
$ws.Range("D8").Value = "pngimages/24_banana.png"
$ws.Range("E8").Value = -0.5
$ws.Range("F8").Value = 0.5

$ws.Range("A9").Value = "trainingaudio/13_kopopi1.wav"
$ws.Range("B9").Value = "pngimages/13_toast.png"
$ws.Range("C9").Value = "trainingaudio/08_tipako2.wav"
$ws.Range("D9").Value = "pngimages/08_bell.png"
$ws.Range("E9").Value = 0.5
$ws.Range("F9").Value = -0.5

$ws.Range("A10").Value = "trainingaudio/11_tokiko1.wav"
$ws.Range("B10").Value = "pngimages/11_compass.png"
$ws.Range("C10").Value = "trainingaudio/10_tokiti1.wav"
$ws.Range("D10").Value = "pngimages/10_backpack.png"
$ws.Range("E10").Value = -0.5
$ws.Range("F10").Value = 0.5

